$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

$ws.Range("B2").Value = 105.98916852820224
$ws.Range("C2").Value = 104.74085481089448
$ws.Range("D2").Value = 107.76606483851549
$ws.Range("E2").Value = 107.10477109939815

$ws.Range("B3").Value = 104.88524901633632
$ws.Range("C3").Value = 107.8360232974745
$ws.Range("D3").Value = 108.64319819792583
$ws.Range("E3").Value = 106.41734465713107

$ws.Range("B1:E3").Select()
